{"js": "// Fix duplicated word \"s\u00fahvezdie\" that appears before \"S\u00fahvezdie Pegasus\"\n// Original: \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Pegasus: ...\"\n// Fixed:    \"V roku 2022 m\u00f4\u017eete pozorova\u0165 S\u00fahvezdie Pegasus: ...\"\n\nconst searchText = \"pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Pegasus\";\nconst replaceText = \"pozorova\u0165 S\u00fahvezdie Pegasus\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Fix duplicated word \"s\u00fahvezdie\" that appears before \"S\u00fahvezdie Pegasus\"\n# Original: \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Pegasus: ...\"\n# Fixed:    \"V roku 2022 m\u00f4\u017eete pozorova\u0165 S\u00fahvezdie Pegasus: ...\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\n    \"pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Pegasus\",  # FindText\n    $true,                                    # MatchCase\n    $false,                                   # MatchWholeWord\n    $false,                                   # MatchWildcards\n    $false,                                   # MatchSoundsLike\n    $false,                                   # MatchAllWordForms\n    $true,                                    # Forward\n    1,                                        # Wrap (wdFindContinue)\n    $false,                                   # Format\n    \"pozorova\u0165 S\u00fahvezdie Pegasus\",            # ReplaceWith\n    2                                         # Replace (wdReplaceAll)\n)\n"}
